# Tutorial 6 solution update:
# - Replace '/' separators in the Date column (A3:A21) with '-'
# - Update a few attendance tally cells (rows 3, 4, 12) to reflect the
#   corrected counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @{
    3  = "28-07-2022"
    4  = "01-08-2022"
    5  = "04-08-2022"
    6  = "08-08-2022"
    7  = "11-08-2022"
    8  = "15-08-2022"
    9  = "18-08-2022"
    10 = "22-08-2022"
    11 = "25-08-2022"
    12 = "29-08-2022"
    13 = "01-09-2022"
    14 = "05-09-2022"
    15 = "08-09-2022"
    16 = "12-09-2022"
    17 = "15-09-2022"
    18 = "19-09-2022"
    19 = "22-09-2022"
    20 = "26-09-2022"
    21 = "29-09-2022"
}

foreach ($row in $dates.Keys) {
    $text = $dates[$row]
    $day = [int]($text.Split("-")[0])
    if ($day -le 12) {
        # Ambiguous as dd-mm-yyyy vs mm-dd-yyyy: force literal text with a
        # leading apostrophe so Excel doesn't silently convert it to a date.
        $ws.Cells.Item($row, 1).Value = "'" + $text
    } else {
        $ws.Cells.Item($row, 1).Value = $text
    }
}

# D = Total Attendance Count, E = Real, F = Duplicate, G = Invalid, H = Absent
$ws.Cells.Item(3, 4).Value = 1   # D3: 0 -> 1
$ws.Cells.Item(3, 7).Value = 1   # G3: 0 -> 1

$ws.Cells.Item(4, 4).Value = 1   # D4: 0 -> 1
$ws.Cells.Item(4, 5).Value = 1   # E4: 0 -> 1
$ws.Cells.Item(4, 8).Value = 0   # H4: 1 -> 0

$ws.Cells.Item(12, 4).Value = 1  # D12: 0 -> 1
$ws.Cells.Item(12, 5).Value = 1  # E12: 0 -> 1
$ws.Cells.Item(12, 8).Value = 0  # H12: 1 -> 0
